$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Error Codes")

# Renumber the error codes in column A (rows 82-97) from the old
# 70001/71001.. series to the new 9xxxx series.
$values = @{
    82 = 90004
    83 = 90005
    84 = 90006
    85 = 90007
    86 = 90008
    87 = 90009
    88 = 90010
    89 = 90011
    90 = 90012
    91 = 90013
    92 = 90014
    93 = 90015
    94 = 90016
    95 = 90017
    96 = 90018
    97 = 90019
}

foreach ($r in $values.Keys) {
    $ws.Range("A$r").Value = $values[$r]
}

# The table uses a repeating 3-row visual pattern down column A: two rows in
# the "normal" style (as seen on A82) followed by one row in the "header of
# group" style (as seen on A84). Re-apply that pattern across A82:A97 so
# every row's look matches its position in the cycle - this folds A85/A86
# (which used to carry a unique right-aligned style) and A88/A89/A91/A92
# (which used to carry the "group header" style) back into the normal
# pattern, and promotes A96 to the group-header look.
$normalSource = $ws.Range("A82")
$groupHeaderSource = $ws.Range("A84")

$normalSource.Copy() | Out-Null
$normalRows = @(82, 83, 85, 86, 88, 89, 91, 92, 94, 95, 97)
foreach ($r in $normalRows) {
    $ws.Range("A$r").PasteSpecial(-4122) | Out-Null
}

$groupHeaderSource.Copy() | Out-Null
$groupHeaderRows = @(84, 87, 90, 93, 96)
foreach ($r in $groupHeaderRows) {
    $ws.Range("A$r").PasteSpecial(-4122) | Out-Null
}

$excel.CutCopyMode = 0

# Leave the selection on the block of rows that was just edited, matching
# where the user would naturally end up after updating this table.
$ws.Range("A82:A97").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 39
